# This script updates the mass-flow result values in the workbook with
# freshly recomputed numbers coming from new input files (see commit
# message: "added new input files"). Only numeric values in the existing
# cells change; layout/structure stays the same.

$wb = $excel.ActiveWorkbook

$wsOutput = $wb.Worksheets.Item("Output_flows")
$wsInput  = $wb.Worksheets.Item("Input_flows")

# ---------------------------------------------------------------------
# Output_flows sheet (columns C, D, E, F)
# ---------------------------------------------------------------------

# Row 7
$wsOutput.Range("C7").Value = [double]"9.528878723240796E-09"
$wsOutput.Range("E7").Value = [double]"5.938818946236801E-10"
$wsOutput.Range("F7").Value = [double]"1.603481115483936E-10"

# Row 12
$wsOutput.Range("C12").Value = [double]"7.617726528932581E-10"
$wsOutput.Range("E12").Value = [double]"2.373852157806095E-11"
$wsOutput.Range("F12").Value = [double]"6.409400826076458E-12"

# Row 13
$wsOutput.Range("C13").Value = [double]"1.644466729824091E-10"
$wsOutput.Range("D13").Value = [double]"1.218731504806904E-11"
$wsOutput.Range("E13").Value = [double]"5.12452223141207E-12"
$wsOutput.Range("F13").Value = [double]"1.383621002481259E-12"

# Row 14
$wsOutput.Range("C14").Value = [double]"8.81093415823931E-11"
$wsOutput.Range("D14").Value = [double]"2.611950208703635E-11"
$wsOutput.Range("E14").Value = [double]"2.745682059389261E-12"
$wsOutput.Range("F14").Value = [double]"7.413341560351006E-13"

# Row 17
$wsOutput.Range("C17").Value = 0.662757391121933
$wsOutput.Range("E17").Value = 0.08261194764822795
$wsOutput.Range("F17").Value = 0.02230522586502155

# Row 18
$wsOutput.Range("C18").Value = 0.1469737584242817
$wsOutput.Range("D18").Value = 0.000435700541392089
$wsOutput.Range("E18").Value = 0.01832011019304685
$wsOutput.Range("F18").Value = 0.004946429752122648

# Row 19
$wsOutput.Range("C19").Value = 0.100597921835844
$wsOutput.Range("D19").Value = 0.001192866994129875
$wsOutput.Range("E19").Value = 0.01253941542342499
$wsOutput.Range("F19").Value = 0.003385642164324747

# ---------------------------------------------------------------------
# Input_flows sheet (column C only)
# ---------------------------------------------------------------------

$wsInput.Range("C7").Value = [double]"1.028310872941287E-08"
$wsInput.Range("C12").Value = [double]"7.666735092058083E-10"
$wsInput.Range("C13").Value = [double]"1.700823802208533E-10"
$wsInput.Range("C14").Value = [double]"1.177158598848538E-10"
$wsInput.Range("C17").Value = 0.7666424305967253
$wsInput.Range("C18").Value = 0.1700795654137782
$wsInput.Range("C19").Value = 0.1177158464177236
